# Update the three-digit x one-digit multiplication problems in the table.
$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "283×8="; New = "574×7=" },
    @{ Old = "496×9="; New = "459×7=" },
    @{ Old = "568×9="; New = "385×7=" },
    @{ Old = "851×6="; New = "134×8=" },
    @{ Old = "996×5="; New = "445×6=" },
    @{ Old = "843×9="; New = "198×6=" },
    @{ Old = "371×9="; New = "237×7=" },
    @{ Old = "895×8="; New = "822×6=" },
    @{ Old = "693×2="; New = "599×8=" },
    @{ Old = "297×7="; New = "806×6=" },
    @{ Old = "942×3="; New = "740×5=" },
    @{ Old = "796×6="; New = "649×6=" },
    @{ Old = "854×6="; New = "719×6=" },
    @{ Old = "446×3="; New = "303×5=" },
    @{ Old = "736×8="; New = "296×8=" },
    @{ Old = "227×7="; New = "319×7=" },
    @{ Old = "282×5="; New = "696×5=" },
    @{ Old = "708×2="; New = "728×2=" },
    @{ Old = "654×3="; New = "548×8=" },
    @{ Old = "923×8="; New = "343×8=" },
    @{ Old = "922×3="; New = "604×3=" },
    @{ Old = "357×2="; New = "788×7=" },
    @{ Old = "222×6="; New = "894×9=" },
    @{ Old = "481×4="; New = "356×5=" },
    @{ Old = "521×7="; New = "527×7=" }
)

foreach ($pair in $replacements) {
    $d.Content.Find.Execute($pair.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $pair.New, 2)
}
